# Insert a new data row at row 2 (pushing the existing rows down by one),
# then fill it in with the new phone/DDD/date values, matching the
# formatting used by the other data rows (style used by row 3, i.e. the
# former row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2; this shifts rows 2..17 down to 3..18.
$ws.Rows.Item(2).Insert()

# Make sure the phone number / DDD / date are written as plain text so
# leading "+" and zeros are preserved (not turned into a number/date).
$ws.Range("A2:C2").NumberFormat = "@"
$ws.Range("A2").Value = "+5513997050892"
$ws.Range("B2").Value = "13"
$ws.Range("C2").Value = "2024-10-18"

# Re-apply the same formatting as the data rows (row 3 here, the old row 2)
# so the new row matches the rest of the table instead of inheriting the
# header's style from the Insert() above.
$ws.Range("A3:C3").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)
